$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 633.3333
$ws.Range("I2").Value = 550
$ws.Range("J2").Value = 800
$ws.Range("K2").Value = 550
$ws.Range("L2").Value = 800
$ws.Range("M2").Value = -437
$ws.Range("N2").Value = -1026

$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws.Range("H64").Value = 2989.35
$ws.Range("I64").Value = 2922
$ws.Range("J64").Value = 3006.1875
$ws.Range("K64").Value = 2922
$ws.Range("L64").Value = 3006.1875
$ws.Range("M64").Value = -2674
$ws.Range("N64").Value = -3502.1875

$ws.Range("H67").Value = 2989.35
$ws.Range("I67").Value = 2922
$ws.Range("J67").Value = 3006.1875
$ws.Range("K67").Value = 2922
$ws.Range("L67").Value = 3006.1875
$ws.Range("M67").Value = -2064
$ws.Range("N67").Value = -4722.1875

$ws.Range("H70").Value = 1388
$ws.Range("I70").Value = 1322.5
$ws.Range("J70").Value = 1431.6666
$ws.Range("K70").Value = 3967.5
$ws.Range("L70").Value = 4294.9998
$ws.Range("M70").Value = -3697.5
$ws.Range("N70").Value = -4834.9998

$ws.Range("H73").Value = 1388
$ws.Range("I73").Value = 1322.5
$ws.Range("J73").Value = 1431.6666
$ws.Range("K73").Value = 3967.5
$ws.Range("L73").Value = 4294.9998
$ws.Range("M73").Value = -3031.5
$ws.Range("N73").Value = -6166.9998

$ws.Range("H76").Value = 1987282.9
$ws.Range("I76").Value = 3231.7144
$ws.Range("K76").Value = 3231.7144
$ws.Range("M76").Value = -2916.7144

$ws.Range("H79").Value = 1987282.9
$ws.Range("I79").Value = 3231.7144
$ws.Range("K79").Value = 3231.7144
$ws.Range("M79").Value = -2139.7144

$ws.Range("H100").Value = 90910580
$ws.Range("I100").Value = 111112150
$ws.Range("K100").Value = 111112150
$ws.Range("M100").Value = -111111609

$ws.Range("H116").Value = 7000.6
$ws.Range("I116").Value = 4000
$ws.Range("J116").Value = 7750.75
$ws.Range("K116").Value = 4000
$ws.Range("L116").Value = 7750.75
$ws.Range("M116").Value = -558
$ws.Range("N116").Value = -14634.75

$ws.Range("H129").Value = 1256.5667
$ws.Range("J129").Value = 1479.875
$ws.Range("L129").Value = 4439.625
$ws.Range("N129").Value = -14439.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 33022.766
$ws.Range("I132").Value = 4116.75
$ws.Range("J132").Value = 58717
$ws.Range("K132").Value = 12350.25
$ws.Range("L132").Value = 176151
$ws.Range("M132").Value = -9820.25
$ws.Range("N132").Value = -181211

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1181.4286
$ws.Range("I99").Value = 1265.6666
$ws.Range("J99").Value = 1069.1111
$ws.Range("K99").Value = 1265.6666
$ws.Range("L99").Value = 1069.1111
$ws.Range("M99").Value = 232.3334
$ws.Range("N99").Value = -4065.1111

$ws.Range("H134").Value = 4699
$ws.Range("I134").Value = 5035.737
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 15107.211
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -12572.211
$ws.Range("N134").Value = -9570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 39119.6
$ws.Range("J133").Value = 39119.6
$ws.Range("L133").Value = 39119.6
$ws.Range("N133").Value = -44179.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 40.8
$ws.Range("I2").Value = 27
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 162
$ws.Range("L2").Value = 300
$ws.Range("M2").Value = -49
$ws.Range("N2").Value = -526

$ws.Range("H5").Value = 1138.6
$ws.Range("I5").Value = 882.36365
$ws.Range("J5").Value = 1339.9286
$ws.Range("K5").Value = 2647.09095
$ws.Range("L5").Value = 4019.7858
$ws.Range("M5").Value = -2535.09095
$ws.Range("N5").Value = -4243.7858

$ws.Range("H26").Value = 554.8333
$ws.Range("J26").Value = 933
$ws.Range("L26").Value = 2799
$ws.Range("N26").Value = -3375

$ws.Range("H40").Value = 102.5
$ws.Range("I40").Value = 63
$ws.Range("K40").Value = 252
$ws.Range("M40").Value = -183

$ws.Range("H86").Value = 41667550
$ws.Range("I86").Value = 728.8570999999999
$ws.Range("J86").Value = 100001100
$ws.Range("K86").Value = 2186.5713
$ws.Range("L86").Value = 300003300
$ws.Range("M86").Value = -1000.5713
$ws.Range("N86").Value = -300005672

$ws.Range("H89").Value = 41667550
$ws.Range("I89").Value = 728.8570999999999
$ws.Range("J89").Value = 100001100
$ws.Range("K89").Value = 6559.7139
$ws.Range("L89").Value = 900009900
$ws.Range("M89").Value = -631.7138999999997
$ws.Range("N89").Value = -900021756

$ws.Range("H113").Value = 436
$ws.Range("I113").Value = 420
$ws.Range("J113").Value = 446.66666
$ws.Range("K113").Value = 1260
$ws.Range("L113").Value = 1339.99998
$ws.Range("M113").Value = 910
$ws.Range("N113").Value = -5679.999980000001

$ws.Range("H117").Value = 1357.25
$ws.Range("I117").Value = 1357.25
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 4071.75
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = -629.75
$ws.Range("N117").ClearContents()

$ws.Range("H123").Value = 4305.3335
$ws.Range("J123").Value = 4960.6
$ws.Range("L123").Value = 14881.8
$ws.Range("N123").Value = -19781.8

$ws.Range("H129").Value = 6634.1904
$ws.Range("I129").Value = 774.44446
$ws.Range("J129").Value = 11029
$ws.Range("K129").Value = 2323.33338
$ws.Range("L129").Value = 33087
$ws.Range("M129").Value = 2676.66662
$ws.Range("N129").Value = -43087

$ws.Range("H131").Value = 700.49
$ws.Range("J131").Value = 718.0538
$ws.Range("L131").Value = 2154.1614
$ws.Range("N131").Value = -12234.1614

$ws.Range("H135").Value = 1138.6
$ws.Range("I135").Value = 882.36365
$ws.Range("J135").Value = 1339.9286
$ws.Range("K135").Value = 7941.27285
$ws.Range("L135").Value = 12059.3574
$ws.Range("M135").Value = -5406.27285
$ws.Range("N135").Value = -17129.3574

$ws.Range("H139").Value = 2041.5927
$ws.Range("I139").Value = 1316.85
$ws.Range("K139").Value = 3950.55
$ws.Range("M139").Value = 1189.45

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5379.3335
$ws.Range("I126").Value = 4378.5713
$ws.Range("K126").Value = 13135.7139
$ws.Range("M126").Value = -10665.7139

$ws.Range("H132").Value = 32262.766
$ws.Range("I132").Value = 4773.75
$ws.Range("J132").Value = 40720.92
$ws.Range("K132").Value = 14321.25
$ws.Range("L132").Value = 122162.76
$ws.Range("M132").Value = -11791.25
$ws.Range("N132").Value = -127222.76

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3160.5293
$ws.Range("I40").Value = 2143.9
$ws.Range("J40").Value = 4612.857
$ws.Range("K40").Value = 2143.9
$ws.Range("L40").Value = 4612.857
$ws.Range("M40").Value = -2007.9
$ws.Range("N40").Value = -4884.857

$ws.Range("H68").Value = 2276.3333
$ws.Range("J68").Value = 2864.5
$ws.Range("L68").Value = 2864.5
$ws.Range("N68").Value = -4362.5

$ws.Range("H71").Value = 2276.3333
$ws.Range("J71").Value = 2864.5
$ws.Range("L71").Value = 14322.5
$ws.Range("N71").Value = -21810.5

$ws.Range("H127").Value = 29882.6
$ws.Range("J127").Value = 29882.6
$ws.Range("L127").Value = 29882.6
$ws.Range("N127").Value = -39802.6

$ws.Range("H132").Value = 805552.1
$ws.Range("I132").Value = 1005774.06
$ws.Range("J132").Value = 4664.3335
$ws.Range("K132").Value = 3017322.18
$ws.Range("L132").Value = 13993.0005
$ws.Range("M132").Value = -3014792.18
$ws.Range("N132").Value = -19053.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 19666.666
$ws.Range("J69").Value = 19666.666
$ws.Range("L69").Value = 19666.666
$ws.Range("N69").Value = -21164.666

$ws.Range("H72").Value = 19666.666
$ws.Range("J72").Value = 19666.666
$ws.Range("L72").Value = 58999.99800000001
$ws.Range("N72").Value = -66487.99800000001
